# Fill in the accuracy / precision / recall values for the new models
# (rows 21-24, columns C:H) in the second ("V2") results table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = 0.91176470588235203
$ws.Range("D21").Value = 0.875
$ws.Range("E21").Value = 0.93333333333333302
$ws.Range("F21").Value = 0.91176470588235203
$ws.Range("G21").Value = 0.9375
$ws.Range("H21").Value = 0.88235294117647001

$ws.Range("C22").Value = 0.91176470588235203
$ws.Range("D22").Value = 0.875
$ws.Range("E22").Value = 0.93333333333333302
$ws.Range("F22").Value = 0.91176470588235203
$ws.Range("G22").Value = 0.9375
$ws.Range("H22").Value = 0.88235294117647001

$ws.Range("C23").Value = 0.82352941176470495
$ws.Range("D23").Value = 0.76470588235294101
$ws.Range("E23").Value = 0.86666666666666603
$ws.Range("F23").Value = 0.88235294117647001
$ws.Range("G23").Value = 0.93333333333333302
$ws.Range("H23").Value = 0.82352941176470495

$ws.Range("C24").Value = 0.79411764705882304
$ws.Range("D24").Value = 0.75
$ws.Range("E24").Value = 0.8
$ws.Range("F24").Value = 0.88235294117647001
$ws.Range("G24").Value = 0.88235294117647001
$ws.Range("H24").Value = 0.88235294117647001

# Match the author's final cursor position from the edit session.
$ws.Range("G30").Select()
